$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.104.03"
$ws.Range("E2").Value = "  -2.56%  "

$ws.Range("D3").Value = "3.043.76"
$ws.Range("E3").Value = "  -5.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").Value = "3.039.78"
$ws.Range("E8").Value = "  -5.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -11.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000213"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -12.56%  "

$ws.Range("D15").Value = "3.522.76"
$ws.Range("E15").Value = "  -5.52%  "

$ws.Range("D16").Value = "64.940.86"
$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("E17").Value = "  -3.84%  "

$ws.Range("D18").Value = "3.025.98"
$ws.Range("E18").Value = "  -5.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "484.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.654"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -13.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.92%  "

$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -14.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.00%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.26%  "

$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "536.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.22%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -12.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0410"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("E39").Value = "  -4.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0782"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -13.53%  "

$ws.Range("D42").Value = "2.788.29"
$ws.Range("E42").Value = "  -4.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.237"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.12%  "

$ws.Range("D48").Value = "0.0₃0507"
$ws.Range("E48").Value = "  -12.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -11.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -16.59%  "
